
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Actual Result" (H) / "Test Status" (I) columns are reset to the
# "not executed" state. Row 7 is intentionally left alone - it keeps its
# PASSED text/green highlighting.
$rowsToReset = @(2, 3, 4, 5, 6, 8, 9, 10, 11, 12, 13, 14, 15, 16)

# D2 carries the plain (no-fill) body style used throughout the sheet - copy
# its formatting onto every "Test Status" cell that must lose its old
# PASSED/FAILED color highlighting.
$ws.Range("D2").Copy()

foreach ($r in $rowsToReset) {
    $ws.Range("I$r").PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

foreach ($r in $rowsToReset) {
    $ws.Range("H$r").Value = "Test not executed"
    $ws.Range("I$r").Value = "Not Run"
}

# Update the test objective wording for TC_011 (row 12).
$ws.Range("C12").Value = "Verify email required validation is visible & Next button in disabled."
